$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cluster names and active case counts (rows 2-13)
$names = @(
    "Billboards The Venue Melbourne",
    "Costa Mushroom Farm Mernda",
    "Crown Melbourne Southbank",
    "Flanagans Border Inn Bacchus Marsh",
    "Melbourne Cricket Ground (MCG)",
    "Monash Health Dandenong Hospital Dandenong",
    "Ms Collins Melbourne",
    "Rupert On Rupert Collingwood",
    "St Vincents Hospital Melbourne Emergency Department Fitzroy",
    "St. Vincent's Hospital Melbourne Fitzroy",
    "The Hatter and the Hare Bayswater",
    "Western Health Sunshine Hospital Emergency Department St Albans"
)

$values = @(10, 21, 11, 14, 20, 16, 40, 11, 25, 14, 14, 12)

# Remove the old rows beyond row 13 (rows 14-22) so the sheet shrinks back down
$oldLastRow = 22
$newLastRow = 13
if ($oldLastRow -gt $newLastRow) {
    $ws.Range("A" + ($newLastRow + 1) + ":B" + $oldLastRow).ClearContents() | Out-Null
}

# Write the new cluster names and active case values into rows 2-13
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
